$d = $word.ActiveDocument

# 1. Split the "Version" run into "Versi" + "on" without altering formatting.
#    A transient bookmark forces a clean run boundary; deleting it afterwards
#    leaves the two runs split but otherwise untouched.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("tmp_split_mark", $splitPoint)
$d.Bookmarks("tmp_split_mark").Delete()

# 2. "Version 1." -> "Version 2" : change the digit, then drop the trailing
#    period from its original run (it gets re-added after the bookmark below).
$digit = $d.Range(8, 9)
$digit.Text = "2"

$oldDot = $d.Range(9, 10)
$oldDot.Text = ""

# 3. Re-insert the period as its own run after the _GoBack bookmark.
$afterBookmark = $d.Range(9, 9)
$afterBookmark.InsertAfter(".")
